$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force text format, assign, then restore default style so no stray formatting remains.
$textCells = @{
    "D5" = "322.63"
    "D6" = "105.15"
    "D7" = "0.524"
    "D9" = "0.543"
    "D10" = "38.07"
    "D14" = "7.18"
    "D17" = "0.849"
    "D19" = "12.79"
    "D20" = "6.58"
    "D22" = "70.67"
    "D23" = "2.41"
    "D24" = "251.61"
    "D26" = "26.21"
    "D28" = "10.04"
    "D29" = "2.21"
    "D30" = "35.16"
    "D32" = "49.45"
    "D33" = "19.75"
    "D34" = "5.37"
    "D35" = "0.0783"
    "D38" = "4.65"
    "D42" = "121.84"
    "D43" = "21.36"
    "D47" = "2.10"
    "D49" = "9.22"
    "D50" = "5.29"
    "D51" = "79.44"
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}

# Remaining cells: plain text assignment (URLs, names, percent strings, multi-dot numbers)
$plainCells = @{
    "D2" = "47.623.79"
    "E2" = "  +4.78%  "
    "D3" = "2.492.33"
    "E3" = "  +2.66%  "
    "E4" = "  -0.02%  "
    "E5" = "  +1.27%  "
    "E6" = "  +2.24%  "
    "E7" = "  +1.54%  "
    "E8" = "  +0.01%  "
    "E9" = "  +2.44%  "
    "E10" = "  +6.88%  "
    "E11" = "  +1.17%  "
    "E12" = "  +1.20%  "
    "E13" = "  +0.84%  "
    "E14" = "  +1.60%  "
    "D15" = "2.881.40"
    "E15" = "  +2.61%  "
    "D16" = "2.497.07"
    "E16" = "  +2.85%  "
    "E17" = "  +0.31%  "
    "D18" = "47.489.36"
    "E18" = "  +4.69%  "
    "E19" = "  +4.52%  "
    "E20" = "  +3.69%  "
    "E21" = "  +1.64%  "
    "E22" = "  +2.60%  "
    "B23" = "ImmutableX"
    "C23" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "E23" = "  +6.28%  "
    "B24" = "BitcoinCash"
    "C24" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "E24" = "  +2.84%  "
    "E25" = "  +3.49%  "
    "E26" = "  +2.05%  "
    "E27" = "  -0.06%  "
    "B28" = "Cosmos"
    "C28" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "E28" = "  +4.63%  "
    "B29" = "Toncoin"
    "C29" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "E29" = "  +1.18%  "
    "E30" = "  +6.35%  "
    "E31" = "  +8.40%  "
    "E32" = "  +0.51%  "
    "E33" = "  -2.95%  "
    "E34" = "  +3.15%  "
    "E35" = "  +1.82%  "
    "E36" = "  +0.06%  "
    "E37" = "  +5.50%  "
    "E38" = "  +4.36%  "
    "E39" = "  +4.37%  "
    "E40" = "  +1.95%  "
    "E41" = "  +1.45%  "
    "E42" = "  -3.44%  "
    "E43" = "  +3.71%  "
    "E44" = "  +2.30%  "
    "D45" = "1.967.07"
    "E45" = "  +2.10%  "
    "E46" = "  +1.77%  "
    "E47" = "  -0.38%  "
    "E48" = "  +0.96%  "
    "E49" = "  +0.91%  "
    "E50" = "  +11.84%  "
    "E51" = "  +3.52%  "
}
foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}